$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "47.340.98"
$ws.Cells.Item(2, 5).Value = "  +3.57%  "

$ws.Cells.Item(3, 4).Value = "2.504.40"
$ws.Cells.Item(3, 5).Value = "  +2.82%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "324.26"
$ws.Cells.Item(5, 5).Value = "  +0.83%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "109.34"
$ws.Cells.Item(6, 5).Value = "  +4.82%  "

$ws.Cells.Item(7, 5).Value = "  +1.87%  "

$ws.Cells.Item(8, 5).Value = "  -0.02%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.543"
$ws.Cells.Item(9, 5).Value = "  +1.54%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "38.94"
$ws.Cells.Item(10, 5).Value = "  +8.40%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "18.59"
$ws.Cells.Item(13, 5).Value = "  +1.74%  "

$ws.Cells.Item(14, 5).Value = "  +2.81%  "

$ws.Cells.Item(15, 4).Value = "2.894.24"
$ws.Cells.Item(15, 5).Value = "  +2.78%  "

$ws.Cells.Item(16, 4).Value = "2.496.65"
$ws.Cells.Item(16, 5).Value = "  +2.39%  "

$ws.Cells.Item(17, 5).Value = "  +2.32%  "

$ws.Cells.Item(18, 4).Value = "47.303.26"
$ws.Cells.Item(18, 5).Value = "  +3.74%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.05"
$ws.Cells.Item(19, 5).Value = "  +5.15%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.74"
$ws.Cells.Item(20, 5).Value = "  +4.95%  "

$ws.Cells.Item(21, 5).Value = "  +2.05%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "71.08"
$ws.Cells.Item(22, 5).Value = "  -0.82%  "

$ws.Cells.Item(23, 5).Value = "  +7.88%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "250.16"
$ws.Cells.Item(24, 5).Value = "  +1.49%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.61"
$ws.Cells.Item(25, 5).Value = "  +3.93%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.19"
$ws.Cells.Item(26, 5).Value = "  +1.74%  "

$ws.Cells.Item(27, 5).Value = "  -0.04%  "

$ws.Cells.Item(28, 5).Value = "  +0.52%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.03"
$ws.Cells.Item(29, 5).Value = "  +3.88%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "35.76"
$ws.Cells.Item(30, 5).Value = "  +6.78%  "

$ws.Cells.Item(31, 5).Value = "  +5.61%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "49.79"
$ws.Cells.Item(32, 5).Value = "  +0.85%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "19.99"
$ws.Cells.Item(33, 5).Value = "  -0.60%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.47"
$ws.Cells.Item(34, 5).Value = "  +4.14%  "

$ws.Cells.Item(35, 5).Value = "  +4.55%  "

$ws.Cells.Item(36, 5).Value = "  +0.21%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.77"
$ws.Cells.Item(37, 5).Value = "  +5.01%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.99"
$ws.Cells.Item(38, 5).Value = "  +5.49%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.00"
$ws.Cells.Item(39, 5).Value = "  +3.16%  "

$ws.Cells.Item(40, 5).Value = "  +1.70%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "122.42"
$ws.Cells.Item(41, 5).Value = "  -3.94%  "

$ws.Cells.Item(42, 5).Value = "  -1.80%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "21.53"
$ws.Cells.Item(43, 5).Value = "  +3.04%  "

$ws.Cells.Item(44, 5).Value = "  +2.76%  "

$ws.Cells.Item(45, 4).Value = "1.990.49"
$ws.Cells.Item(45, 5).Value = "  +1.62%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.06"
$ws.Cells.Item(46, 5).Value = "  +3.37%  "

$ws.Cells.Item(47, 5).Value = "  -1.10%  "

$ws.Cells.Item(48, 5).Value = "  -1.23%  "

$ws.Cells.Item(49, 5).Value = "  -0.80%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.34"
$ws.Cells.Item(50, 5).Value = "  +10.16%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "78.55"
$ws.Cells.Item(51, 5).Value = "  +1.75%  "
